# Add the four "resulting categories" labels below the existing header row.
# These become new shared strings and new rows (21-24) in column A, styled
# with a slightly larger (12pt) plain font - matching the look of the
# existing "Stories for the Sprint" header column but without the bold/
# fill treatment used on row 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21").Value = "Remaining Effort"
$ws.Range("A22").Value = "Ideal Burndown"
$ws.Range("A23").Value = "Stories Remaining"
$ws.Range("A24").Value = "Stories Completed"

# Give the new category labels their own (4th) font/style: 12pt, no bold,
# no fill - distinct from the existing 3 cell styles already in the sheet.
$ws.Range("A21:A24").Font.Size = 12
